$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# so numeric-looking strings (e.g. "3.80", "1.40") are not coerced into numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '57.712.69'
$ws.Range('E2').Value = '  +2.68%  '
$ws.Range('D3').Value = '2.443.63'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '158.24'
$ws.Range('E5').Value = '  +6.09%  '
$ws.Range('D6').Value = '493.25'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').Value = '0.994'
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +21.68%  '
$ws.Range('D9').Value = '2.452.88'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('E10').Value = '  +12.46%  '
$ws.Range('D11').Value = '0.102'
$ws.Range('E11').Value = '  +4.28%  '
$ws.Range('D12').Value = '0.335'
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '2.856.01'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '57.598.86'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '20.97'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '2.447.33'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').Value = '4.71'
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('D20').Value = '330.64'
$ws.Range('E20').Value = '  +4.87%  '
$ws.Range('D21').Value = '10.12'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '5.95'
$ws.Range('E23').Value = '  +4.20%  '
$ws.Range('D24').Value = '58.52'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('D25').Value = '0.406'
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').Value = '0.994'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('D28').Value = '2.527.11'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Value = '7.37'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').Value = '0.0₃0807'
$ws.Range('E30').Value = '  +3.46%  '
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '18.87'
$ws.Range('E32').Value = '  +5.54%  '
$ws.Range('D33').Value = '150.75'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  +3.53%  '
$ws.Range('D35').Value = '5.38'
$ws.Range('E35').Value = '  +6.75%  '
$ws.Range('D36').Value = '1.16'
$ws.Range('E36').Value = '  +3.45%  '
$ws.Range('D37').Value = '3.80'
$ws.Range('E37').Value = '  +4.84%  '
$ws.Range('D38').Value = '0.831'
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.40'
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '34.32'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('D41').Value = '3.57'
$ws.Range('E41').Value = '  +4.78%  '
$ws.Range('D42').Value = '0.101'
$ws.Range('E42').Value = '  +6.64%  '
$ws.Range('D43').Value = '279.45'
$ws.Range('E43').Value = '  +7.54%  '
$ws.Range('D44').Value = '0.989'
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('D45').Value = '0.602'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').Value = '0.0541'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('D47').Value = '0.0231'
$ws.Range('E47').Value = '  +3.61%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').Value = '4.70'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('D50').Value = '18.07'
$ws.Range('E50').Value = '  +5.23%  '
$ws.Range('D51').Value = '0.692'
$ws.Range('E51').Value = '  +11.12%  '
